# Apply the "make forex gains from dividends tax-free" edit to the
# "Foreign Currencies" sheet (and its knock-on effect on the
# "ELSTER - Summary" sheet, which sums the "Gains (incl. losses)" row).
#
# Rows 3 and 4 of "Foreign Currencies" describe USD that came from
# dividend payments (not bought), so they should not be taxed as FOREX
# gains: their Gain [EUR] (column G) becomes 0 and their Comment
# (column H) is replaced with a new explanatory note instead of the old
# "Held for N days ..." text.

$wb = $excel.ActiveWorkbook
$fc = $wb.Worksheets.Item("Foreign Currencies")
$elster = $wb.Worksheets.Item("ELSTER - Summary")

$newComment = "FOREX not acquired (e.g. received dividend payments), thus gains not taxed."

# Row 2 (unrelated currency lot) - quantity correction.
$fc.Range("B2").Value = 1247.91

# Row 3 - USD received as dividend payment, no taxable forex gain.
$fc.Range("G3").Value = 0
$fc.Range("H3").Value = $newComment

# Row 4 - USD received as dividend payment, no taxable forex gain.
$fc.Range("G4").Value = 0
$fc.Range("H4").Value = $newComment

# Row 5 - recomputed buy quantity/gain.
$fc.Range("B5").Value = 2567.09
$fc.Range("G5").Value = 57.06

# Row 6 - recomputed buy quantity/gain.
$fc.Range("B6").Value = 849.87
$fc.Range("G6").Value = -43.76

# Row 7 - recomputed buy quantity/gain.
$fc.Range("B7").Value = 135.13
$fc.Range("G7").Value = -7.45

# Summary rows 9-11 (Gains incl. losses / Gains excl. losses / Losses).
$fc.Range("G9").Value = 5.85
$fc.Range("G10").Value = 57.06
$fc.Range("G11").Value = -51.21

# The ELSTER - Summary sheet mirrors the "Gains (incl. losses)" total for
# foreign currencies (row 7, "Zeilen 42 - 48: Gewinn / Verlust aus Verkauf
# von Fremdwährungen") - keep it in sync.
$elster.Range("C7").Value = 5.85
